$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# --- Bold restyle of existing "Range model" title and "contrasts" header (reuses existing bold style) ---
$ws2.Range("C2").Font.Bold = $true
$ws2.Range("C7").Font.Bold = $true

# --- "tilts" header: bold + right aligned ---
$ws2.Range("B7").Font.Bold = $true
$ws2.Range("B7").HorizontalAlignment = -4152   # xlRight

# --- Rows 8-13: add per-tilt observation text in column C, with wrap + 45pt row height ---
# Values are entered in the same order the shared strings were first authored so the
# resulting shared-string table ordering matches the source workbook.
$ws2.Range("C8").Value = "High Conf: 0.33`nTf Acc: 0.52`nMy Acc: 0.55"
$ws2.Range("C10").Value = "High Conf: 0.35`nTf Acc: 0.56`nMy Acc: 0.61"
$ws2.Range("C11").Value = "High Conf: 0.36`nTf Acc: 0.61`nMy Acc: 0.67"
$ws2.Range("C12").Value = "High Conf: 0.43`nTf Acc: 0.68`nMy Acc: 0.75"
$ws2.Range("C13").Value = "High Conf: 0.48`nTf Acc: 0.72`nMy Acc: 0.82"
$ws2.Range("C9").Value = "High Conf: 0.33`nTf Acc: 0.53`nMy Acc: 0.56"

$ws2.Range("C8:C13").WrapText = $true
$ws2.Rows.Item(8).RowHeight = 45
$ws2.Rows.Item(9).RowHeight = 45
$ws2.Rows.Item(10).RowHeight = 45
$ws2.Rows.Item(11).RowHeight = 45
$ws2.Rows.Item(12).RowHeight = 45
$ws2.Rows.Item(13).RowHeight = 45

# --- Row 16-20: Observations section ---
$ws2.Range("C16").Value = "Observations"
$ws2.Range("C16").WrapText = $true

$ws2.Range("D16").Value = "The High Conf (0.38) and My Acc (0.66) calculated using the test set comprised of the full range of tilts and contrasts reflected an average of all the High Conf and My Acc calculations across tilts"

$ws2.Range("D19").Value = "The confidence began slowly ticking up at the third tilt (0.4) and took a large jump in the fifth tilt (1.6)"

$ws2.Range("D17").Value = "The first two tilts (0.1 & 0.2) were too slight for the model to discriminate between very well."

$ws2.Range("D18").Value = "Not until the 4th tilt (0.8) was there a more noticeable performance above chance (0.6+)"

$ws2.Range("D20").Value = "The model becomes more confident and performs better as the tilts become greater"

# --- sheetView: remove topLeftCell scroll position, update selection to D21 ---
$ws2.Activate()
$ws2.Range("D21").Select()

Write-Host "edit complete"
